$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stale per-row date formatting rows (6,7,8) and the second data
# row (4) -- the whole used range collapses down to just two data rows.
$ws.Rows("4:8").Delete()

# The D/E columns carried manual width overrides (customWidth/bestFit) --
# clear that column-level formatting so no <cols> survive in the sheet.
$ws.Columns("D:E").ClearFormats()

# Row 1 (new first record)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 21312312
$ws.Range("C1").Value = "выавыавыава"
$ws.Range("D1").Value = "26.01.2024"
$ws.Range("E1").Value = "31.01.2024"

# Row 2 (existing record, content replaced)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 13213213
$ws.Range("C2").Value = "вавыаывыва"
$ws.Range("D2").Value = "26.01.2024"
$ws.Range("E2").Value = "30.01.2024"

$ws.Range("D1:E2").NumberFormat = "dd.mm.yyyy"

$ws.Range("I8").Select()
